$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1980830670926517
$ws.Range("C2").Value = 0.549520766773163
$ws.Range("J2").Value = 0.01277955271565495
$ws.Range("P2").Value = 0.1309904153354633
$ws.Range("S2").Value = 0.1086261980830671
$ws.Range("B3").Value = 0.005555555555555556
$ws.Range("C3").Value = 0.03888888888888889
$ws.Range("J3").Value = 0.01666666666666667
$ws.Range("P3").Value = 0.7611111111111111
$ws.Range("S3").Value = 0.1777777777777778
$ws.Range("J4").Value = 0.0625
$ws.Range("P4").Value = 0.625
$ws.Range("S4").Value = 0.3125
$ws.Range("B6").Value = 0.04639175257731959
$ws.Range("E6").Value = 0.005154639175257732
$ws.Range("F6").Value = 0.03092783505154639
$ws.Range("J6").Value = 0.3092783505154639
$ws.Range("O6").Value = 0.02061855670103093
$ws.Range("Q6").Value = 0.1649484536082474
$ws.Range("R6").Value = 0.07216494845360824
$ws.Range("S6").Value = 0.3505154639175257
$ws.Range("B7").Value = 0.1092896174863388
$ws.Range("D7").Value = 0.02185792349726776
$ws.Range("E7").Value = 0.00546448087431694
$ws.Range("F7").Value = 0.04918032786885246
$ws.Range("J7").Value = 0.1639344262295082
$ws.Range("O7").Value = 0.01639344262295082
$ws.Range("Q7").Value = 0.1366120218579235
$ws.Range("R7").Value = 0.07650273224043716
$ws.Range("S7").Value = 0.4207650273224044
$ws.Range("B8").Value = 0.09895833333333333
$ws.Range("D8").Value = 0.015625
$ws.Range("E8").Value = 0.002604166666666667
$ws.Range("F8").Value = 0.04427083333333334
$ws.Range("J8").Value = 0.09114583333333333
$ws.Range("O8").Value = 0.03385416666666666
$ws.Range("Q8").Value = 0.15625
$ws.Range("R8").Value = 0.08854166666666667
$ws.Range("S8").Value = 0.46875
$ws.Range("B9").Value = 0.1176470588235294
$ws.Range("D9").Value = 0.005882352941176471
$ws.Range("E9").Value = 0.01176470588235294
$ws.Range("F9").Value = 0.06470588235294118
$ws.Range("J9").Value = 0.1176470588235294
$ws.Range("O9").Value = 0.02352941176470588
$ws.Range("Q9").Value = 0.1176470588235294
$ws.Range("R9").Value = 0.06470588235294118
$ws.Range("S9").Value = 0.4764705882352941
$ws.Range("B10").Value = 0.1417525773195876
$ws.Range("D10").Value = 0.03264604810996564
$ws.Range("E10").Value = 0.001718213058419244
$ws.Range("F10").Value = 0.07216494845360824
$ws.Range("J10").Value = 0.09450171821305842
$ws.Range("O10").Value = 0.02663230240549828
$ws.Range("Q10").Value = 0.1898625429553265
$ws.Range("R10").Value = 0.06701030927835051
$ws.Range("S10").Value = 0.3737113402061856
$ws.Range("G11").Value = 0.1678321678321678
$ws.Range("J11").Value = 0.09090909090909091
$ws.Range("K11").Value = 0.1993006993006993
$ws.Range("L11").Value = 0.5314685314685315
$ws.Range("S11").Value = 0.01048951048951049
$ws.Range("G12").Value = 0.7051282051282052
$ws.Range("J12").Value = 0.1987179487179487
$ws.Range("K12").Value = 0.00641025641025641
$ws.Range("L12").Value = 0.02564102564102564
$ws.Range("S12").Value = 0.0641025641025641
$ws.Range("G13").Value = 0.6222222222222222
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("S13").Value = 0.04444444444444445
$ws.Range("F15").Value = 0.01680672268907563
$ws.Range("H15").Value = 0.134453781512605
$ws.Range("I15").Value = 0.05882352941176471
$ws.Range("J15").Value = 0.3571428571428572
$ws.Range("K15").Value = 0.06302521008403361
$ws.Range("M15").Value = 0.004201680672268907
$ws.Range("O15").Value = 0.06722689075630252
$ws.Range("S15").Value = 0.2983193277310924
$ws.Range("F16").Value = 0.0198019801980198
$ws.Range("H16").Value = 0.1237623762376238
$ws.Range("I16").Value = 0.0891089108910891
$ws.Range("J16").Value = 0.405940594059406
$ws.Range("K16").Value = 0.1089108910891089
$ws.Range("M16").Value = 0.0396039603960396
$ws.Range("O16").Value = 0.06930693069306931
$ws.Range("S16").Value = 0.1435643564356436
$ws.Range("F17").Value = 0.02247191011235955
$ws.Range("H17").Value = 0.1685393258426966
$ws.Range("I17").Value = 0.1067415730337079
$ws.Range("J17").Value = 0.3876404494382023
$ws.Range("K17").Value = 0.09550561797752809
$ws.Range("M17").Value = 0.02528089887640449
$ws.Range("N17").Value = 0.002808988764044944
$ws.Range("O17").Value = 0.06460674157303371
$ws.Range("S17").Value = 0.1264044943820225
$ws.Range("H18").Value = 0.1390728476821192
$ws.Range("I18").Value = 0.09271523178807947
$ws.Range("J18").Value = 0.4834437086092715
$ws.Range("K18").Value = 0.1258278145695364
$ws.Range("M18").Value = 0.006622516556291391
$ws.Range("O18").Value = 0.05298013245033113
$ws.Range("S18").Value = 0.09933774834437085
$ws.Range("F19").Value = 0.02289452166802943
$ws.Range("H19").Value = 0.2035977105478332
$ws.Range("I19").Value = 0.07113654946852004
$ws.Range("J19").Value = 0.3802125919869174
$ws.Range("K19").Value = 0.1120196238757155
$ws.Range("M19").Value = 0.02207686017988553
$ws.Range("N19").Value = 0.0008176614881439084
$ws.Range("O19").Value = 0.07686017988552739
$ws.Range("S19").Value = 0.1103843008994276
